$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 109.026058
$ws.Range("H2").Value = 327.078174
$ws.Range("I2").Value = 0.3049840938689738
$ws.Range("J2").Value = 0.3049840938689738
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 223.2367336666667
$ws.Range("N2").Value = 669.710201
$ws.Range("O2").Value = 0.9523995969492647
$ws.Range("P2").Value = 0.9523995969492646
$ws.Range("Q2").Value = 24338.62107247255
$ws.Range("R2").Value = 219047.589652253
$ws.Range("S2").Value = 0.2904667280767473
$ws.Range("T2").Value = 0.2904667280767473

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 109.026058
$ws.Range("H3").Value = 327.078174
$ws.Range("I3").Value = 0.3049840938689738
$ws.Range("J3").Value = 0.3049840938689738
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.280784333333333
$ws.Range("N3").Value = 12.842353
$ws.Range("O3").Value = 0.01826320071400582
$ws.Range("P3").Value = 0.01826320071400582
$ws.Range("Q3").Value = 466.7170410114913
$ws.Range("R3").Value = 4200.453369103421
$ws.Range("S3").Value = 0.005569985720908259
$ws.Range("T3").Value = 0.005569985720908259

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 109.026058
$ws.Range("H4").Value = 327.078174
$ws.Range("I4").Value = 0.3049840938689738
$ws.Range("J4").Value = 0.3049840938689738
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.876463666666666
$ws.Range("N4").Value = 20.629391
$ws.Range("O4").Value = 0.02933720233672951
$ws.Range("P4").Value = 0.0293372023367295
$ws.Range("Q4").Value = 749.7137265568925
$ws.Range("R4").Value = 6747.423539012033
$ws.Range("S4").Value = 0.008947380071318189
$ws.Range("T4").Value = 0.008947380071318187

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 89.97721833333333
$ws.Range("H5").Value = 269.931655
$ws.Range("I5").Value = 0.2516978134001918
$ws.Range("J5").Value = 0.2516978134001917
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 223.2367336666667
$ws.Range("N5").Value = 669.710201
$ws.Range("O5").Value = 0.9523995969492647
$ws.Range("P5").Value = 0.9523995969492646
$ws.Range("Q5").Value = 20086.22032514585
$ws.Range("R5").Value = 180775.9829263126
$ws.Range("S5").Value = 0.2397168960353539
$ws.Range("T5").Value = 0.2397168960353538

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 89.97721833333333
$ws.Range("H6").Value = 269.931655
$ws.Range("I6").Value = 0.2516978134001918
$ws.Range("J6").Value = 0.2516978134001917
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.280784333333333
$ws.Range("N6").Value = 12.842353
$ws.Range("O6").Value = 0.01826320071400582
$ws.Range("P6").Value = 0.01826320071400582
$ws.Range("Q6").Value = 385.1730665982461
$ws.Range("R6").Value = 3466.557599384214
$ws.Range("S6").Value = 0.004596807685404085
$ws.Range("T6").Value = 0.004596807685404084

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 89.97721833333333
$ws.Range("H7").Value = 269.931655
$ws.Range("I7").Value = 0.2516978134001918
$ws.Range("J7").Value = 0.2516978134001917
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.876463666666666
$ws.Range("N7").Value = 20.629391
$ws.Range("O7").Value = 0.02933720233672951
$ws.Range("P7").Value = 0.0293372023367295
$ws.Range("Q7").Value = 618.7250726969005
$ws.Range("R7").Value = 5568.525654272104
$ws.Range("S7").Value = 0.007384109679433813
$ws.Range("T7").Value = 0.00738410967943381

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 158.477852
$ws.Range("H8").Value = 475.433556
$ws.Range("I8").Value = 0.4433180927308344
$ws.Range("J8").Value = 0.4433180927308344
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 223.2367336666667
$ws.Range("N8").Value = 669.710201
$ws.Range("O8").Value = 0.9523995969492647
$ws.Range("P8").Value = 0.9523995969492646
$ws.Range("Q8").Value = 35378.07803898941
$ws.Range("R8").Value = 318402.7023509047
$ws.Range("S8").Value = 0.4222159728371634
$ws.Range("T8").Value = 0.4222159728371634

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 158.477852
$ws.Range("H9").Value = 475.433556
$ws.Range("I9").Value = 0.4433180927308344
$ws.Range("J9").Value = 0.4433180927308344
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.280784333333333
$ws.Range("N9").Value = 12.842353
$ws.Range("O9").Value = 0.01826320071400582
$ws.Range("P9").Value = 0.01826320071400582
$ws.Range("Q9").Value = 678.4095060219186
$ws.Range("R9").Value = 6105.685554197267
$ws.Range("S9").Value = 0.008096407307693473
$ws.Range("T9").Value = 0.008096407307693471

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 158.477852
$ws.Range("H10").Value = 475.433556
$ws.Range("I10").Value = 0.4433180927308344
$ws.Range("J10").Value = 0.4433180927308344
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.876463666666666
$ws.Range("N10").Value = 20.629391
$ws.Range("O10").Value = 0.02933720233672951
$ws.Range("P10").Value = 0.0293372023367295
$ws.Range("Q10").Value = 1089.767191249377
$ws.Range("R10").Value = 9807.904721244395
$ws.Range("S10").Value = 0.0130057125859775
$ws.Range("T10").Value = 0.0130057125859775
